$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column F: re-measured "DistanceThreshold" values -------------------
# The original column E ("DistanceThreshold") is being relabeled to
# "old73to75" and a fresh column F ("DistanceThreshold") is added with the
# updated mean values (bw 0.73 to 0.76).

# Copy the plain header/row formatting (style used by column A) onto the new
# column F so it matches the rest of the non-highlighted columns.
$ws.Range("A1:A13").Copy()
$ws.Range("F1:F13").PasteSpecial(-4122)

# Headers
$ws.Range("E1").Value = "old73to75"
$ws.Range("F1").Value = "DistanceThreshold"

# New column F data (rows 2-13)
$ws.Range("F2").Value = 149
$ws.Range("F3").Value = 134
$ws.Range("F4").Value = 166
$ws.Range("F5").Value = 158
$ws.Range("F6").Value = 172
$ws.Range("F7").Value = 165
$ws.Range("F8").Value = 160
$ws.Range("F9").Value = 118
$ws.Range("F10").Value = 175
$ws.Range("F11").Value = 145
$ws.Range("F12").Value = 246
$ws.Range("F13").Value = 205

# Corrected values in the (renamed) original column E
$ws.Range("E6").Value = 170
$ws.Range("E12").Value = 246

# --- Remove the unused filler block (rows 14-20, partial rows 15 & 21) ------
$ws.Range("J14:Q14").Clear()
$ws.Range("J15").Clear()
$ws.Range("Q15").Clear()
$ws.Range("J16:Q20").Clear()
$ws.Range("J21:K21").Clear()

# --- Cosmetic: match the author's last selection ----------------------------
$null = $ws.Range("E13").Select()
